$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text to uppercase / new labels
$ws.Range("A1").Value = "NOMBRE"
$ws.Range("B1").Value = "DEPTO."
$ws.Range("C1").Value = "CORREO"
$ws.Range("D1").Value = "DEUDAS"

# New column E: header + style matching the other header cells
$ws.Range("E1").Value = "FECHA DE MODIFICACIÓN"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# New data cell for the existing data row
$ws.Range("E2").Value = "2020/12/09, 18:55:05"

# Set width of the new column E (21.6 ColumnWidth units renders as 22.5 in the saved XML)
$ws.Columns.Item(5).ColumnWidth = 21.6

# Update the selection to match the target state
$ws.Range("G7").Select()

$wb.Save()
